# Apply the CERScore recalculation described in the commit:
# "FIX - CERScore como inverso da distancia. Recalculados valores."
#
# The method rows are reordered and their C:F values (CERScore distances)
# are recalculated. Two methods (revise, cem-vae) are removed entirely,
# shrinking the table from 13 data rows to 11 data rows (dimension A1:F12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data definitions in final desired order (row 1 = header, row 2 = KAOGExp stay as-is)
$data = @(
    @("face-knn",     1,    15.53369195207143, 49.13011562692634, 1128.518752337495,  108.3578934703077),
    @("dice",         1,    33.98934500833334, 59.35752915839514, 105.7838445710516,  87.16566050438655),
    @("cchvae",       1,    9.628893003238096, 20.95153073975084, 29.99150357596599,  55.26166336831749),
    @("clue",         1,    19.4032325785,     50.82135965033636, 89.80743963893373,  111.7907486823867),
    @("cruds",        1,    17.57648285559524, 50.66177961118662, 85.92389996374665,  99.95343782682214),
    @("cem",          1,    20.63110641733333, 100.8232904601979, 349.8987473563467,  125.4288327353597),
    @("wachter",      0.77, 6.2926589535,      159.3113481599378, 11304.82193388034,  601.2263835808923),
    @("face-epsilon", 1,    15.89281101821429, 52.83096651546239, 1086.996173822396,  109.0872926566667),
    @("gs",           1,    13.71322545833333, 266.278842631892,  28006.83932611461,  459.8404036686949),
    @("ar",           0.53, 10.20333166333333, 14.8592423095449,  29.66194059493817,  18.37380823208987)
)

# First, remove the two trailing rows (13 and 14) so the remaining rows 3..12
# line up exactly with the new data (delete from the bottom up so row
# indices of earlier rows stay valid).
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()

# Now overwrite rows 3..12 with the recalculated values, in the new order.
$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
